$d = $word.ActiveDocument

# Replace "rasi bintang Orion" with "Rasi bintang Orion" throughout the document
# (capitalizes the word "rasi" -> "Rasi" wherever it precedes "bintang Orion")
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("rasi bintang Orion", $true, $false, $false, $false, $false, `
               $true, 1, $false, "Rasi bintang Orion", 2)
